# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# across the four sheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 19
$ws1.Range("F5").Value = 19831
$ws1.Range("F7").Value = 2288
$ws1.Range("F9").Value = 612
$ws1.Range("F11").Value = 696
$ws1.Range("F16").Value = 77
$ws1.Range("F19").Value = 197

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 199
$ws2.Range("F4").Value = 13
$ws2.Range("F5").Value = 17
$ws2.Range("G6").Value = "不可售"
$ws2.Range("F7").Value = 291
$ws2.Range("F8").Value = 133
$ws2.Range("F10").Value = 15
$ws2.Range("F16").Value = 84

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6007
$ws3.Range("F3").Value = 650
$ws3.Range("F4").Value = 595

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6007
$ws4.Range("F3").Value = 650
$ws4.Range("F4").Value = 595
$ws4.Range("F5").Value = 199
$ws4.Range("F6").Value = 19
$ws4.Range("F10").Value = 19831
$ws4.Range("F11").Value = 13
$ws4.Range("F12").Value = 17
$ws4.Range("G14").Value = "不可售"
$ws4.Range("F15").Value = 291
$ws4.Range("F16").Value = 2288
$ws4.Range("F18").Value = 133
$ws4.Range("F19").Value = 612
$ws4.Range("F21").Value = 696
$ws4.Range("F26").Value = 15
$ws4.Range("F29").Value = 77
$ws4.Range("F36").Value = 197
$ws4.Range("F37").Value = 84
$ws4.Range("F38").Value = 84
